$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1375
$ws1.Range("F7").Value = 599
$ws1.Range("F10").Value = 26
$ws1.Range("F11").Value = 656

# Sheet "全部类型" (All types) - same events, update corresponding rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 1375
$ws4.Range("F20").Value = 599
$ws4.Range("F24").Value = 26
$ws4.Range("F26").Value = 656
